$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet: update row 2 (91bbbb67 file) datetimes
$wsZh.Range("E2").Value = "2016-03-24 10:21:09"
$wsZh.Range("H2").Value = "2016-03-24 10:22:06"

# de-de sheet: update row 2 (91bbbb67 file) datetimes
$wsDe.Range("E2").Value = "2016-03-24 10:21:18"
$wsDe.Range("H2").Value = "2016-03-24 10:22:21"
